# feat: add 2022-Q1 data
# 1) Insert a new "2022-Q1" sheet right before the "总计" (summary) sheet,
#    using the existing "2021-Q4" sheet as a formatting template so headers /
#    index-column styling match the other quarterly sheets.
$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("A1:H18").Copy()
$q1.Range("A1").PasteSpecial(-4122)
$q1.Range("A1").ClearContents()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Columns B:G hold text values (fund code/name keep leading zeros & fixed
# decimal formatting), so force Text format before assignment; only column A
# (row index) and H (rank) are real numbers.
$q1.Range("B2:G18").NumberFormat = "@"
$q1.Cells.Item(2,1).Value = 0
$q1.Cells.Item(2,2).Value = "001305"
$q1.Cells.Item(2,3).Value = "九泰天富改革新动力混合A"
$q1.Cells.Item(2,4).Value = "3.74"
$q1.Cells.Item(2,5).Value = "88.86"
$q1.Cells.Item(2,6).Value = "6.53"
$q1.Cells.Item(2,7).Value = "0.2442"
$q1.Cells.Item(2,8).Value = 5
$q1.Cells.Item(3,1).Value = 1
$q1.Cells.Item(3,2).Value = "006013"
$q1.Cells.Item(3,3).Value = "易方达鑫转招利混合A"
$q1.Cells.Item(3,4).Value = "12.76"
$q1.Cells.Item(3,5).Value = "25.24"
$q1.Cells.Item(3,6).Value = "1.35"
$q1.Cells.Item(3,7).Value = "0.1723"
$q1.Cells.Item(3,8).Value = 5
$q1.Cells.Item(4,1).Value = 2
$q1.Cells.Item(4,2).Value = "001782"
$q1.Cells.Item(4,3).Value = "九泰久益灵活配置混合A"
$q1.Cells.Item(4,4).Value = "2.33"
$q1.Cells.Item(4,5).Value = "94.33"
$q1.Cells.Item(4,6).Value = "7.04"
$q1.Cells.Item(4,7).Value = "0.1640"
$q1.Cells.Item(4,8).Value = 6
$q1.Cells.Item(5,1).Value = 3
$q1.Cells.Item(5,2).Value = "001543"
$q1.Cells.Item(5,3).Value = "宝盈新锐灵活配置混合A"
$q1.Cells.Item(5,4).Value = "3.21"
$q1.Cells.Item(5,5).Value = "93.26"
$q1.Cells.Item(5,6).Value = "4.31"
$q1.Cells.Item(5,7).Value = "0.1384"
$q1.Cells.Item(5,8).Value = 10
$q1.Cells.Item(6,1).Value = 4
$q1.Cells.Item(6,2).Value = "001844"
$q1.Cells.Item(6,3).Value = "九泰久益灵活配置混合C"
$q1.Cells.Item(6,4).Value = "1.47"
$q1.Cells.Item(6,5).Value = "94.33"
$q1.Cells.Item(6,6).Value = "7.04"
$q1.Cells.Item(6,7).Value = "0.1035"
$q1.Cells.Item(6,8).Value = 6
$q1.Cells.Item(7,1).Value = 5
$q1.Cells.Item(7,2).Value = "004128"
$q1.Cells.Item(7,3).Value = "新疆前海联合泳隆灵活配置混合A"
$q1.Cells.Item(7,4).Value = "0.86"
$q1.Cells.Item(7,5).Value = "91.05"
$q1.Cells.Item(7,6).Value = "4.91"
$q1.Cells.Item(7,7).Value = "0.0422"
$q1.Cells.Item(7,8).Value = 6
$q1.Cells.Item(8,1).Value = 6
$q1.Cells.Item(8,2).Value = "007040"
$q1.Cells.Item(8,3).Value = "新疆前海联合泳隆灵活配置混合C"
$q1.Cells.Item(8,4).Value = "0.82"
$q1.Cells.Item(8,5).Value = "91.05"
$q1.Cells.Item(8,6).Value = "4.91"
$q1.Cells.Item(8,7).Value = "0.0403"
$q1.Cells.Item(8,8).Value = 6
$q1.Cells.Item(9,1).Value = 7
$q1.Cells.Item(9,2).Value = "009912"
$q1.Cells.Item(9,3).Value = "九泰天富改革新动力混合C"
$q1.Cells.Item(9,4).Value = "0.59"
$q1.Cells.Item(9,5).Value = "88.86"
$q1.Cells.Item(9,6).Value = "6.53"
$q1.Cells.Item(9,7).Value = "0.0385"
$q1.Cells.Item(9,8).Value = 5
$q1.Cells.Item(10,1).Value = 8
$q1.Cells.Item(10,2).Value = "006014"
$q1.Cells.Item(10,3).Value = "易方达鑫转招利混合C"
$q1.Cells.Item(10,4).Value = "1.88"
$q1.Cells.Item(10,5).Value = "25.24"
$q1.Cells.Item(10,6).Value = "1.35"
$q1.Cells.Item(10,7).Value = "0.0254"
$q1.Cells.Item(10,8).Value = 5
$q1.Cells.Item(11,1).Value = 9
$q1.Cells.Item(11,2).Value = "000066"
$q1.Cells.Item(11,3).Value = "诺安鸿鑫混合"
$q1.Cells.Item(11,4).Value = "0.74"
$q1.Cells.Item(11,5).Value = "81.34"
$q1.Cells.Item(11,6).Value = "3.08"
$q1.Cells.Item(11,7).Value = "0.0228"
$q1.Cells.Item(11,8).Value = 9
$q1.Cells.Item(12,1).Value = 10
$q1.Cells.Item(12,2).Value = "007578"
$q1.Cells.Item(12,3).Value = "宝盈新锐灵活配置混合C"
$q1.Cells.Item(12,4).Value = "0.20"
$q1.Cells.Item(12,5).Value = "93.26"
$q1.Cells.Item(12,6).Value = "4.31"
$q1.Cells.Item(12,7).Value = "0.0086"
$q1.Cells.Item(12,8).Value = 10
$q1.Cells.Item(13,1).Value = 11
$q1.Cells.Item(13,2).Value = "009336"
$q1.Cells.Item(13,3).Value = "平安中证500指数增强A"
$q1.Cells.Item(13,4).Value = "0.30"
$q1.Cells.Item(13,5).Value = "87.90"
$q1.Cells.Item(13,6).Value = "2.26"
$q1.Cells.Item(13,7).Value = "0.0068"
$q1.Cells.Item(13,8).Value = 6
$q1.Cells.Item(14,1).Value = 12
$q1.Cells.Item(14,2).Value = "008437"
$q1.Cells.Item(14,3).Value = "九泰行业优选灵活配置混合A"
$q1.Cells.Item(14,4).Value = "0.11"
$q1.Cells.Item(14,5).Value = "51.13"
$q1.Cells.Item(14,6).Value = "5.41"
$q1.Cells.Item(14,7).Value = "0.0060"
$q1.Cells.Item(14,8).Value = 5
$q1.Cells.Item(15,1).Value = 13
$q1.Cells.Item(15,2).Value = "000892"
$q1.Cells.Item(15,3).Value = "九泰天宝灵活配置混合A"
$q1.Cells.Item(15,4).Value = "0.07"
$q1.Cells.Item(15,5).Value = "90.81"
$q1.Cells.Item(15,6).Value = "4.63"
$q1.Cells.Item(15,7).Value = "0.0032"
$q1.Cells.Item(15,8).Value = 7
$q1.Cells.Item(16,1).Value = 14
$q1.Cells.Item(16,2).Value = "008438"
$q1.Cells.Item(16,3).Value = "九泰行业优选灵活配置混合C"
$q1.Cells.Item(16,4).Value = "0.06"
$q1.Cells.Item(16,5).Value = "51.13"
$q1.Cells.Item(16,6).Value = "5.41"
$q1.Cells.Item(16,7).Value = "0.0032"
$q1.Cells.Item(16,8).Value = 5
$q1.Cells.Item(17,1).Value = 15
$q1.Cells.Item(17,2).Value = "009337"
$q1.Cells.Item(17,3).Value = "平安中证500指数增强C"
$q1.Cells.Item(17,4).Value = "0.14"
$q1.Cells.Item(17,5).Value = "87.90"
$q1.Cells.Item(17,6).Value = "2.26"
$q1.Cells.Item(17,7).Value = "0.0032"
$q1.Cells.Item(17,8).Value = 6
$q1.Cells.Item(18,1).Value = 16
$q1.Cells.Item(18,2).Value = "002028"
$q1.Cells.Item(18,3).Value = "九泰天宝灵活配置混合C"
$q1.Cells.Item(18,4).Value = "0.00"
$q1.Cells.Item(18,5).Value = "90.81"
$q1.Cells.Item(18,6).Value = "4.63"
$q1.Range("G18").NumberFormat = "General"
$q1.Cells.Item(18,7).Value = 0
$q1.Cells.Item(18,8).Value = 7
$q1.Range("B2:G18").ClearFormats()

# 2) Update the "总计" (summary) sheet - it is now the last sheet - inserting
#    a new top data row for 2022-Q1 and bumping the existing row index values.
$total = $wb.Worksheets.Item($wb.Worksheets.Count)
$total.Rows.Item(2).Insert(-4121)
$total.Range("B2:D2").ClearFormats()
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 17
$total.Cells.Item(2,4).Value = 1.02

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(7,1).Value = 5
